# doc: modificacion objetivo general
#
# Updates the "Objetivo General" slide (slide 5) so that the goal
# paragraph reads "...optimizar la organización del centro de formación,
# mejorar la comunicación..." instead of "...optimizar la organización
# institucional, mejorar la comunicación...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text.IndexOf("institucional, mejorar la comunicación") -ge 0) {
            $needle = "institucional, "
            $fullText = $tr.Text
            $pos = $fullText.IndexOf($needle)
            if ($pos -ge 0) {
                $target = $tr.Characters($pos + 1, $needle.Length)
                $target.Text = "del centro de formación, "
            }
        }
    }
}
